# "Little bits of tidying during recent debug"
#
# 1. constants sheet: remove the blank spacer row (old row 8), shifting the
#    econ_* rows up by one.
# 2. time_variants sheet: drop the unused scenario_7/8/9 columns (BI:BK),
#    and replace the BB17 formula with its static computed value.
# 3. Break the orphaned external link to data_fiji.xlsx.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("constants")
$ws1.Rows.Item(8).Delete()

$ws2 = $wb.Worksheets.Item("time_variants")
$ws2.Range("BB17").Value = 50
$ws2.Range("BI1:BK1").EntireColumn.Delete()

$links = $wb.LinkSources()
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}
